# Updates cryptos list price/volume figures (and reorders a couple of
# coin rows) per the latest GitHub Actions data refresh.
# Note: some "Price" values look like plain numbers (e.g. "96.00"); a
# leading apostrophe forces Excel to keep them as literal text, matching
# the original inline-string cell content (e.g. "96.00" instead of 96).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.664.14'
$ws.Range("E2").Value = '  -0.33%  '
$ws.Range("D3").Value = '2.290.10'
$ws.Range("E3").Value = '  -1.06%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''96.00'
$ws.Range("E5").Value = '  +1.93%  '
$ws.Range("D6").Value = '''268.36'
$ws.Range("E6").Value = '  -0.77%  '
$ws.Range("D7").Value = '''0.622'
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = '''0.609'
$ws.Range("E9").Value = '  -2.40%  '
$ws.Range("D10").Value = '''45.42'
$ws.Range("E10").Value = '  +1.09%  '
$ws.Range("D11").Value = '''0.0933'
$ws.Range("E11").Value = '  -0.60%  '
$ws.Range("D12").Value = '''7.91'
$ws.Range("E12").Value = '  -3.27%  '
$ws.Range("E13").Value = '  +0.89%  '
$ws.Range("D14").Value = '2.631.75'
$ws.Range("E14").Value = '  -1.19%  '
$ws.Range("D15").Value = '''15.34'
$ws.Range("E15").Value = '  -0.26%  '
$ws.Range("D16").Value = '''0.849'
$ws.Range("E16").Value = '  -1.59%  '
$ws.Range("D17").Value = '2.285.57'
$ws.Range("E17").Value = '  -1.22%  '
$ws.Range("D18").Value = '43.571.34'
$ws.Range("E18").Value = '  -0.58%  '
$ws.Range("D19").Value = '''0.0000109'
$ws.Range("E19").Value = '  +2.06%  '
$ws.Range("D20").Value = '''6.20'
$ws.Range("E20").Value = '  -1.56%  '
$ws.Range("D21").Value = '''72.08'
$ws.Range("E21").Value = '  +0.63%  '
$ws.Range("E22").Value = '  +11.88%  '
$ws.Range("D23").Value = '''232.59'
$ws.Range("E23").Value = '  -2.78%  '
$ws.Range("D24").Value = '''9.12'
$ws.Range("E24").Value = '  -5.47%  '
$ws.Range("D25").Value = '''2.60'
$ws.Range("E25").Value = '  +3.22%  '
$ws.Range("E26").Value = '  -0.15%  '
$ws.Range("D27").Value = '''11.22'
$ws.Range("E27").Value = '  -1.67%  '
$ws.Range("E28").Value = '  +2.24%  '
$ws.Range("D29").Value = '''40.07'
$ws.Range("E29").Value = '  +2.53%  '
$ws.Range("D30").Value = '''2.28'
$ws.Range("E30").Value = '  -3.69%  '
$ws.Range("D31").Value = '''174.89'
$ws.Range("E31").Value = '  +1.71%  '
$ws.Range("D32").Value = '''21.80'
$ws.Range("E32").Value = '  -3.35%  '
$ws.Range("D33").Value = '''0.0893'
$ws.Range("E33").Value = '  -0.78%  '
$ws.Range("D34").Value = '''5.36'
$ws.Range("E34").Value = '  -3.67%  '
$ws.Range("E35").Value = '  -0.86%  '
$ws.Range("E36").Value = '  -4.03%  '
$ws.Range("D37").Value = '''0.0351'
$ws.Range("E37").Value = '  -3.07%  '
$ws.Range("D38").Value = '''4.37'
$ws.Range("E38").Value = '  -3.19%  '
$ws.Range("E39").Value = '  -4.07%  '
$ws.Range("D40").Value = '''0.239'
$ws.Range("E40").Value = '  +1.69%  '
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").Value = '''12.34'
$ws.Range("E42").Value = '  +0.59%  '
$ws.Range("D43").Value = '''65.43'
$ws.Range("E43").Value = '  +5.40%  '
$ws.Range("D44").Value = '''1.35'
$ws.Range("E44").Value = '  +0.71%  '
$ws.Range("D45").Value = '''8.78'
$ws.Range("E45").Value = '  -2.44%  '
$ws.Range("E46").Value = '  -0.37%  '
$ws.Range("D47").Value = '''5.15'
$ws.Range("E47").Value = '  -6.03%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = '''96.72'
$ws.Range("E48").Value = '  -3.80%  '
$ws.Range("B49").Value = 'TrustWalletToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D49").Value = '''1.19'
$ws.Range("E49").Value = '  -1.79%  '
$ws.Range("D50").Value = '''0.433'
$ws.Range("E50").Value = '  +0.77%  '
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D51").Value = '''0.185'
$ws.Range("E51").Value = '  +6.83%  '
